# PresentatieSprint2.pptx — "flightpath reset bij stop knop"
#
# 1) Slide 3 title: "Backlog: sprint 2 (gerealiseerd)"
#       -> "Backlog: sprint 2 (gerealiseerd 💪)"   (💪 rendered in Wingdings)
# 2) Slide 4 title: "Backlog: sprint 2 (niet gerealiseerd)" -> bumped to 48pt
#       and "niet gerealiseerd)" -> "niet gerealiseerd 😢)"
# 3) Slide 6 bullet: "Vliegroutes, obstakels en " + "scan zones" runs
#       merged into a single run "Vliegroutes, obstakels en scan zones"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 3 - add the flexed-biceps glyph (Wingdings) to the title
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange

# Rebuild the trailing ": sprint 2 (gerealiseerd)" run as three runs:
# ": sprint 2 (gerealiseerd " / <symbol> / ")"
$run3b = $title3.Runs(2, 1)
$run3b.Text = ": sprint 2 (gerealiseerd )"

# position of the space just before the closing paren, where the glyph goes
$insertAt = $title3.Characters($title3.Length, 1).Start
$title3.Characters($insertAt, 1).InsertBefore("💪")

# give the new glyph its own run formatted with the Wingdings symbol font
$glyph3 = $title3.Characters($insertAt, 1)
$glyph3.Font.Name = "Wingdings"

# ---------------------------------------------------------------------
# 2) Slide 4 - bump title size to 48pt and add the crying-face glyph
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1).TextFrame.TextRange

$run4b = $title4.Runs(2, 1)
$run4b.Text = ": sprint 2 (niet gerealiseerd 😢)"

$title4.Font.Size = 48

# ---------------------------------------------------------------------
# 3) Slide 6 - merge "Vliegroutes, obstakels en " + "scan zones" into one run
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$body6 = $s6.Shapes.Item(2).TextFrame.TextRange

$target6 = $body6.Paragraphs(9, 1)
$firstRun6 = $body6.Characters($target6.Start, 26)
$firstRun6.Delete()

$merged6 = $body6.Paragraphs(9, 1)
$merged6.Text = "PLACEHOLDER__"
$merged6b = $body6.Paragraphs(9, 1)
$merged6b.Text = "Vliegroutes, obstakels en scan zones"
